$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Move the "R21..R24" block (old cols K:N) left into E:H, overwriting the
#     old "R4..R7" contents, then wipe out everything from I to P (rows 1-5) ---

# Row 1 (headers, shared-string cells)
$ws.Range("E1").Value = "R21"
$ws.Range("F1").Value = "R22"
$ws.Range("G1").Value = "R23"
$ws.Range("H1").Value = "R24"

# Row 2
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 10

# Row 3
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 20

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 4

# Clear out the now-stale tail columns.
$ws.Range("I1:P4").ClearContents() | Out-Null
$ws.Range("I5:P5").ClearContents() | Out-Null

# Column A (rows 2-5) loses its conditional-style formatting.
$ws.Range("A2:A5").ClearFormats() | Out-Null

# Row 6 and the lone A7 marker cell go away entirely.
$ws.Range("A6:A7").Clear() | Out-Null

# Update the active selection shown in the sheet view.
$ws.Range("J6").Select() | Out-Null

# Window position tweak recorded in the workbook view.
$excel.ActiveWindow.Left = 6880
